# The document holds a single results table ("ranovas table") with columns
# Variable | h2 | QST | CVA. Several QST (and one h2) values were corrected
# (values approximately halved) to fix a small calculation error.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    # Exclude the trailing end-of-cell marker so we replace only the
    # visible text, not the cell's paragraph mark.
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText
}

Set-CellValue $t 2  3 "0.174"   # Latex exudation: QST 0.296 -> 0.174
Set-CellValue $t 8  3 "0.116"   # Weevil damage (quantitative): QST 0.209 -> 0.116
Set-CellValue $t 9  3 "0.492"   # Flowering success: QST 0.660 -> 0.492
Set-CellValue $t 11 3 "0.305"   # Flower size: QST 0.467 -> 0.305
Set-CellValue $t 12 3 "0.555"   # Flowering duration: QST 0.714 -> 0.555
Set-CellValue $t 13 2 "0.000"   # Date of first flower: h2 0.039 -> 0.000
Set-CellValue $t 14 3 "0.287"   # Follicles: QST 0.446 -> 0.287
Set-CellValue $t 15 3 "0.002"   # Date of first follicle: QST 0.005 -> 0.002
Set-CellValue $t 17 3 "0.184"   # D. plexippus abundance: QST 0.311 -> 0.184
Set-CellValue $t 18 3 "0.081"   # L. asclepiadis abundance: QST 0.151 -> 0.081
Set-CellValue $t 20 3 "0.575"   # LDMC: QST 0.730 -> 0.575
Set-CellValue $t 21 3 "0.524"   # SLA: QST 0.687 -> 0.524
Set-CellValue $t 22 3 "0.031"   # Height before flowering: QST 0.060 -> 0.031
Set-CellValue $t 23 3 "0.052"   # Height after flowering: QST 0.098 -> 0.052
Set-CellValue $t 25 3 "0.009"   # Ramets before flowering: QST 0.017 -> 0.009
Set-CellValue $t 26 3 "0.017"   # Ramets after flowering: QST 0.034 -> 0.017
Set-CellValue $t 27 3 "0.139"   # Mortality: QST 0.245 -> 0.139

# Row-height auto-layout recalculation artifact picked up alongside the
# "Flower size" row edit above (auto height, 571 -> 572 twips, i.e. 28.55 -> 28.6 pt).
$t.Rows.Item(11).Height = 28.6

Write-Output "done"
